$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value2 = 2836
$ws.Range("I34").Value2 = 2836
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 2836
$ws.Range("L34").Value2 = 0
$ws.Range("M34").Value2 = -2633
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value2 = 2836
$ws.Range("I36").Value2 = 2836
$ws.Range("J36").Value2 = 0
$ws.Range("K36").Value2 = 2836
$ws.Range("L36").Value2 = 0
$ws.Range("M36").Value2 = -2121
$ws.Range("N36").ClearContents()
$ws.Range("H51").Value2 = 2708.5
$ws.Range("I51").Value2 = 2000
$ws.Range("K51").Value2 = 2000
$ws.Range("M51").Value2 = -1516
$ws.Range("H62").Value2 = 3865.5
$ws.Range("I62").Value2 = 3832
$ws.Range("J62").Value2 = 3966
$ws.Range("K62").Value2 = 3832
$ws.Range("L62").Value2 = 3966
$ws.Range("M62").Value2 = -3208
$ws.Range("N62").Value2 = -5214
$ws.Range("H65").Value2 = 3865.5
$ws.Range("I65").Value2 = 3832
$ws.Range("J65").Value2 = 3966
$ws.Range("K65").Value2 = 19160
$ws.Range("L65").Value2 = 19830
$ws.Range("M65").Value2 = -16040
$ws.Range("N65").Value2 = -26070
$ws.Range("H133").Value2 = 70000
$ws.Range("J133").Value2 = 70000
$ws.Range("L133").Value2 = 70000
$ws.Range("N133").Value2 = -80120
$ws.Range("H135").Value2 = 191.58333
$ws.Range("I135").Value2 = 191.58333
$ws.Range("K135").Value2 = 1724.24997
$ws.Range("M135").Value2 = 810.7500300000002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4174.037
$ws.Range("I32").Value2 = 4956.2383
$ws.Range("K32").Value2 = 4956.2383
$ws.Range("M32").Value2 = -4669.2383
$ws.Range("H45").Value2 = 2123.1428
$ws.Range("I45").Value2 = 2123.1428
$ws.Range("K45").Value2 = 2123.1428
$ws.Range("M45").Value2 = -1746.1428
$ws.Range("H122").Value2 = 1559.0714
$ws.Range("I122").Value2 = 1557.8
$ws.Range("J122").Value2 = 1562.25
$ws.Range("K122").Value2 = 4673.4
$ws.Range("L122").Value2 = 4686.75
$ws.Range("M122").Value2 = -2223.4
$ws.Range("N122").Value2 = -9586.75
$ws.Range("H132").Value2 = 1861.925
$ws.Range("I132").Value2 = 1568.8889
$ws.Range("K132").Value2 = 4706.6667
$ws.Range("M132").Value2 = -2176.6667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 50001500
$ws.Range("I99").Value2 = 55557056
$ws.Range("K99").Value2 = 55557056
$ws.Range("M99").Value2 = -55555558
$ws.Range("H104").Value2 = 75000
$ws.Range("J104").Value2 = 75000
$ws.Range("L104").Value2 = 75000
$ws.Range("N104").Value2 = -81988
$ws.Range("H105").Value2 = 333335780
$ws.Range("I105").Value2 = 500002500
$ws.Range("J105").Value2 = 2309
$ws.Range("K105").Value2 = 500002500
$ws.Range("L105").Value2 = 2309
$ws.Range("M105").Value2 = -500000753
$ws.Range("N105").Value2 = -5803
$ws.Range("H106").Value2 = 24235
$ws.Range("J106").Value2 = 24235
$ws.Range("L106").Value2 = 24235
$ws.Range("N106").Value2 = -26759
$ws.Range("H107").Value2 = 1650.4615
$ws.Range("I107").Value2 = 1211.2667
$ws.Range("J107").Value2 = 2249.3635
$ws.Range("K107").Value2 = 1211.2667
$ws.Range("L107").Value2 = 2249.3635
$ws.Range("M107").Value2 = 708.7333000000001
$ws.Range("N107").Value2 = -6089.363499999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 100001050
$ws.Range("I16").Value2 = 111112110
$ws.Range("J16").Value2 = 1500
$ws.Range("K16").Value2 = 111112110
$ws.Range("L16").Value2 = 1500
$ws.Range("M16").Value2 = -111111823
$ws.Range("N16").Value2 = -2074
$ws.Range("H31").Value2 = 992
$ws.Range("I31").Value2 = 992
$ws.Range("K31").Value2 = 992
$ws.Range("M31").Value2 = -697
$ws.Range("H34").Value2 = 992
$ws.Range("I34").Value2 = 992
$ws.Range("K34").Value2 = 992
$ws.Range("M34").Value2 = -790
$ws.Range("H58").Value2 = 522.2963
$ws.Range("I58").Value2 = 452.0909
$ws.Range("J58").Value2 = 831.2
$ws.Range("K58").Value2 = 452.0909
$ws.Range("L58").Value2 = 831.2
$ws.Range("M58").Value2 = -249.0909
$ws.Range("N58").Value2 = -1237.2
$ws.Range("H60").Value2 = 11897.85
$ws.Range("I60").Value2 = 0
$ws.Range("K60").Value2 = 0
$ws.Range("M60").ClearContents()
$ws.Range("H113").Value2 = 100001050
$ws.Range("I113").Value2 = 111112110
$ws.Range("J113").Value2 = 1500
$ws.Range("K113").Value2 = 111112110
$ws.Range("L113").Value2 = 1500
$ws.Range("M113").Value2 = -111109940
$ws.Range("N113").Value2 = -5840
$ws.Range("H136").Value2 = 522.2963
$ws.Range("I136").Value2 = 452.0909
$ws.Range("J136").Value2 = 831.2
$ws.Range("K136").Value2 = 1356.2727
$ws.Range("L136").Value2 = 2493.6
$ws.Range("M136").Value2 = 1193.7273
$ws.Range("N136").Value2 = -7593.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 107.166664
$ws.Range("J2").Value2 = 164.5
$ws.Range("L2").Value2 = 987
$ws.Range("N2").Value2 = -1213
$ws.Range("H7").Value2 = 626.8
$ws.Range("I7").Value2 = 683
$ws.Range("J7").Value2 = 402
$ws.Range("K7").Value2 = 2049
$ws.Range("L7").Value2 = 1206
$ws.Range("M7").Value2 = -1937
$ws.Range("N7").Value2 = -1430
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 6105.6665
$ws.Range("I80").Value2 = 1700
$ws.Range("J80").Value2 = 6986.8
$ws.Range("K80").Value2 = 1700
$ws.Range("L80").Value2 = 6986.8
$ws.Range("M80").Value2 = -702
$ws.Range("N80").Value2 = -8982.799999999999
$ws.Range("H83").Value2 = 6105.6665
$ws.Range("I83").Value2 = 1700
$ws.Range("J83").Value2 = 6986.8
$ws.Range("K83").Value2 = 8500
$ws.Range("L83").Value2 = 34934
$ws.Range("M83").Value2 = -3508
$ws.Range("N83").Value2 = -44918
$ws.Range("H122").Value2 = 1292.6666
$ws.Range("I122").Value2 = 1287.7778
$ws.Range("J122").Value2 = 1300
$ws.Range("K122").Value2 = 3863.3334
$ws.Range("L122").Value2 = 3900
$ws.Range("M122").Value2 = -1413.3334
$ws.Range("N122").Value2 = -8800
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 868.13635
$ws.Range("I16").Value2 = 855
$ws.Range("K16").Value2 = 855
$ws.Range("M16").Value2 = -685
$ws.Range("H22").Value2 = 872.3461
$ws.Range("I22").Value2 = 1088.6154
$ws.Range("K22").Value2 = 1088.6154
$ws.Range("M22").Value2 = -793.6153999999999
$ws.Range("H27").Value2 = 872.3461
$ws.Range("I27").Value2 = 1088.6154
$ws.Range("K27").Value2 = 1088.6154
$ws.Range("M27").Value2 = -981.6153999999999
$ws.Range("H122").Value2 = 27781242
$ws.Range("I122").Value2 = 41669696
$ws.Range("J122").Value2 = 4334.6665
$ws.Range("K122").Value2 = 125009088
$ws.Range("L122").Value2 = 13003.9995
$ws.Range("M122").Value2 = -125006638
$ws.Range("N122").Value2 = -17903.9995
$ws.Range("H132").Value2 = 36402.758
$ws.Range("I132").Value2 = 1658.8
$ws.Range("J132").Value2 = 113611.555
$ws.Range("K132").Value2 = 4976.4
$ws.Range("L132").Value2 = 340834.665
$ws.Range("M132").Value2 = -2446.4
$ws.Range("N132").Value2 = -345894.665
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 580.1818
$ws.Range("J107").Value2 = 564.3333
$ws.Range("L107").Value2 = 1692.9999
$ws.Range("N107").Value2 = -5532.9999
$ws.Range("H122").Value2 = 23644172
$ws.Range("I122").Value2 = 26008288
$ws.Range("K122").Value2 = 78024864
$ws.Range("M122").Value2 = -78022414
